$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update collaborator organization name and remove stray tab before "KICT"
$ws.Range("B2").Value = "`t University of Kansai"
$ws.Range("D2").Value = "KICT"

# Reflect the new active selection on the sheet (no frozen/scrolled top-left cell)
$ws.Range("F13").Select()
